$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell D4 gets the text value "R"
$ws.Range("D4").Value = "R"

# Cell L5 changes value from 4 to 3
$ws.Range("L5").Value = 3

# Update the selection shown in the sheet view to K4
$ws.Range("K4").Select()
